$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price/volume snapshot (GitHub Actions data refresh).
# Rows 8/9 (Cardano <-> OKB) and rows 15/16 (ShibaInu <-> Chainlink) also swap
# position because the source ranking reordered them.
$updates = @(
    @{ Cell = 'D2'; Value = '23.388.41' }
    @{ Cell = 'E2'; Value = '  -0.15%  ' }
    @{ Cell = 'D3'; Value = '1.629.06' }
    @{ Cell = 'E3'; Value = '  -0.72%  ' }
    @{ Cell = 'D4'; Value = '0.9995' }
    @{ Cell = 'E4'; Value = '  -0.06%  ' }
    @{ Cell = 'D5'; Value = '0.9997' }
    @{ Cell = 'E5'; Value = '  -0.04%  ' }
    @{ Cell = 'D6'; Value = '302.12' }
    @{ Cell = 'E6'; Value = '  -0.93%  ' }
    @{ Cell = 'D7'; Value = '0.3765' }
    @{ Cell = 'E7'; Value = '  +0.78%  ' }
    @{ Cell = 'B8'; Value = 'OKB' }
    @{ Cell = 'C8'; Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb' }
    @{ Cell = 'D8'; Value = '51.72' }
    @{ Cell = 'E8'; Value = '  -1.10%  ' }
    @{ Cell = 'B9'; Value = 'Cardano' }
    @{ Cell = 'C9'; Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada' }
    @{ Cell = 'D9'; Value = '0.3630' }
    @{ Cell = 'E9'; Value = '  +0.01%  ' }
    @{ Cell = 'D10'; Value = '0.08175' }
    @{ Cell = 'E10'; Value = '  +0.67%  ' }
    @{ Cell = 'D11'; Value = '1.221' }
    @{ Cell = 'E11'; Value = '  -2.55%  ' }
    @{ Cell = 'D12'; Value = '0.9998' }
    @{ Cell = 'E12'; Value = '  -0.05%  ' }
    @{ Cell = 'D13'; Value = '22.23' }
    @{ Cell = 'E13'; Value = '  -2.43%  ' }
    @{ Cell = 'D14'; Value = '6.470' }
    @{ Cell = 'E14'; Value = '  -1.90%  ' }
    @{ Cell = 'B15'; Value = 'Chainlink' }
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link' }
    @{ Cell = 'D15'; Value = '7.307' }
    @{ Cell = 'E15'; Value = '  +0.34%  ' }
    @{ Cell = 'B16'; Value = 'ShibaInu' }
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib' }
    @{ Cell = 'D16'; Value = '0.00001239' }
    @{ Cell = 'E16'; Value = '  -2.51%  ' }
    @{ Cell = 'D17'; Value = '1.617.22' }
    @{ Cell = 'E17'; Value = '  -0.87%  ' }
    @{ Cell = 'D18'; Value = '94.70' }
    @{ Cell = 'E18'; Value = '  +0.34%  ' }
    @{ Cell = 'D19'; Value = '0.06957' }
    @{ Cell = 'E19'; Value = '  +0.93%  ' }
    @{ Cell = 'D20'; Value = '17.55' }
    @{ Cell = 'E20'; Value = '  -3.11%  ' }
    @{ Cell = 'D21'; Value = '6.540' }
    @{ Cell = 'E21'; Value = '  +0.49%  ' }
    @{ Cell = 'D22'; Value = '0.9999' }
    @{ Cell = 'E22'; Value = '  -0.08%  ' }
    @{ Cell = 'D23'; Value = '12.51' }
    @{ Cell = 'E23'; Value = '  -2.03%  ' }
    @{ Cell = 'D24'; Value = '23.373.63' }
    @{ Cell = 'E24'; Value = '  -0.27%  ' }
    @{ Cell = 'D25'; Value = '2.492' }
    @{ Cell = 'E25'; Value = '  +3.27%  ' }
    @{ Cell = 'D26'; Value = '3.080' }
    @{ Cell = 'E26'; Value = '  +0.01%  ' }
    @{ Cell = 'D27'; Value = '21.16' }
    @{ Cell = 'E27'; Value = '  -0.09%  ' }
    @{ Cell = 'D28'; Value = '150.24' }
    @{ Cell = 'E28'; Value = '  -0.87%  ' }
    @{ Cell = 'D29'; Value = '5.277' }
    @{ Cell = 'E29'; Value = '  -0.99%  ' }
    @{ Cell = 'D30'; Value = '132.84' }
    @{ Cell = 'E30'; Value = '  -2.13%  ' }
    @{ Cell = 'D31'; Value = '1.799.90' }
    @{ Cell = 'E31'; Value = '  -0.63%  ' }
    @{ Cell = 'D32'; Value = '6.614' }
    @{ Cell = 'E32'; Value = '  -2.97%  ' }
    @{ Cell = 'E33'; Value = '  -5.41%  ' }
    @{ Cell = 'D34'; Value = '1.063' }
    @{ Cell = 'E34'; Value = '  +11.90%  ' }
    @{ Cell = 'D35'; Value = '11.30' }
    @{ Cell = 'E35'; Value = '  +9.11%  ' }
    @{ Cell = 'D36'; Value = '0.02762' }
    @{ Cell = 'E36'; Value = '  -1.69%  ' }
    @{ Cell = 'D37'; Value = '0.2491' }
    @{ Cell = 'E37'; Value = '  -1.06%  ' }
    @{ Cell = 'D38'; Value = '0.08760' }
    @{ Cell = 'E38'; Value = '  +0.04%  ' }
    @{ Cell = 'D39'; Value = '0.07143' }
    @{ Cell = 'E39'; Value = '  -1.12%  ' }
    @{ Cell = 'D40'; Value = '5.965' }
    @{ Cell = 'E40'; Value = '  -2.27%  ' }
    @{ Cell = 'D41'; Value = '0.6986' }
    @{ Cell = 'E41'; Value = '  -1.00%  ' }
    @{ Cell = 'D42'; Value = '1.327' }
    @{ Cell = 'E42'; Value = '  -3.39%  ' }
    @{ Cell = 'D43'; Value = '15.72' }
    @{ Cell = 'E43'; Value = '  -1.40%  ' }
    @{ Cell = 'D44'; Value = '12.01' }
    @{ Cell = 'E44'; Value = '  -3.66%  ' }
    @{ Cell = 'D45'; Value = '0.6452' }
    @{ Cell = 'E45'; Value = '  -1.03%  ' }
    @{ Cell = 'D46'; Value = '0.9993' }
    @{ Cell = 'E46'; Value = '  +0.03%  ' }
    @{ Cell = 'D47'; Value = '2.273' }
    @{ Cell = 'E47'; Value = '  -2.29%  ' }
    @{ Cell = 'D48'; Value = '3.960' }
    @{ Cell = 'E48'; Value = '  -1.34%  ' }
    @{ Cell = 'D49'; Value = '0.07978' }
    @{ Cell = 'E49'; Value = '  +0.10%  ' }
    @{ Cell = 'D50'; Value = '126.29' }
    @{ Cell = 'E50'; Value = '  -1.71%  ' }
    @{ Cell = 'D51'; Value = '1.188' }
    @{ Cell = 'E51'; Value = '  -0.93%  ' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    # Force Text format so numeric-looking strings (e.g. "0.9995", "23.388.41")
    # are stored as text rather than being auto-converted to numbers.
    $range.NumberFormat = "@"
    $range.Value = $u.Value
}
